$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new meeting diary entry on row 14, copying formatting from row 13
# (the most recently added entry) so the new row matches the sheet's look.
$ws.Range("A13:E13").Copy($ws.Range("A14:E14"))
$ws.Rows.Item(14).RowHeight = $ws.Rows.Item(13).RowHeight

$ws.Range("A14").Value = Get-Date -Year 2023 -Month 10 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("B14").Value = 0.70833333333333337
$ws.Range("C14").Value = 0.8125
$ws.Range("D14").Value = "All"
$ws.Range("E14").Value = "Finished working on Regression and Created the presentation"

# Column E needs to grow a bit to fit the new (slightly longer) discussion text.
$ws.Columns.Item(5).ColumnWidth = 68

# Clear the clipboard marching-ants state left by Copy.
$excel.CutCopyMode = 0

# Match the workbook's final selection state.
$ws.Range("D17").Select()

$wb.Save()
